$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.250.68'
$ws.Range("E2").Value = '  +0.31%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.906.35'
$ws.Range("E3").Value = '  +0.02%  '

$ws.Range("E4").Value = '  +0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '307.59'
$ws.Range("E5").Value = '  +0.38%  '

$ws.Range("E6").Value = '  +0.09%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5252'
$ws.Range("E7").Value = '  +0.31%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3809'
$ws.Range("E8").Value = '  +1.04%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07290'
$ws.Range("E9").Value = '  +0.53%  '

$ws.Range("E10").Value = '  +2.15%  '

$ws.Range("E11").Value = '  -0.02%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08188'
$ws.Range("E12").Value = '  -3.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '96.40'
$ws.Range("E13").Value = '  -0.52%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.368'
$ws.Range("E14").Value = '  +1.39%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.539.84'
$ws.Range("E15").Value = '  -19.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.000'
$ws.Range("E16").Value = '  +0.03%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008675'
$ws.Range("E17").Value = '  +0.06%  '

$ws.Range("E18").Value = '  +1.45%  '

$ws.Range("E19").Value = '  +0.08%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '27.283.79'

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.117'
$ws.Range("E21").Value = '  +0.59%  '

$ws.Range("E22").Value = '  +1.80%  '

$ws.Range("E23").Value = '  +1.07%  '

$ws.Range("E24").Value = '  +0.96%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '150.13'
$ws.Range("E25").Value = '  +2.10%  '

$ws.Range("E26").Value = '  +0.09%  '

$ws.Range("E27").Value = '  -0.52%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '116.70'
$ws.Range("E28").Value = '  +1.29%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.848'
$ws.Range("E29").Value = '  +0.56%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.854'
$ws.Range("E30").Value = '  -1.54%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09249'
$ws.Range("E31").Value = '  -0.60%  '

$ws.Range("E32").Value = '  +3.77%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05069'
$ws.Range("E33").Value = '  +0.21%  '

$ws.Range("E34").Value = '  -1.27%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.985'
$ws.Range("E35").Value = '  +1.14%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.354'
$ws.Range("E36").Value = '  -2.65%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.737'
$ws.Range("E37").Value = '  +4.92%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.5774'
$ws.Range("E38").Value = '  +0.69%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02004'
$ws.Range("E39").Value = '  +0.24%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.078'
$ws.Range("E40").Value = '  +0.33%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.113'
$ws.Range("E41").Value = '  -0.28%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.632'
$ws.Range("E42").Value = '  -0.11%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '117.30'
$ws.Range("E43").Value = '  +1.14%  '

$ws.Range("E44").Value = '  +0.45%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4916'
$ws.Range("E45").Value = '  +0.98%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.19'
$ws.Range("E46").Value = '  -0.25%  '

$ws.Range("E47").Value = '  +0.09%  '

$ws.Range("E48").Value = '  +1.38%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '38.86'
$ws.Range("E49").Value = '  +3.01%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '64.47'
$ws.Range("E50").Value = '  +0.38%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.06055'
$ws.Range("E51").Value = '  +1.70%  '
